$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-7 with new TPM values
# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cfh"
$ws.Range("C2").Value = "Itgam"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.4128076666666667
$ws.Range("H2").Value = 1.238423
$ws.Range("I2").Value = 0.001366259689176221
$ws.Range("J2").Value = 0.001366259689176221
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.005673666666666667
$ws.Range("N2").Value = 0.017021
$ws.Range("O2").Value = 0.1234828534325781
$ws.Range("P2").Value = 0.1234828534325781
$ws.Range("Q2").Value = 0.002342133098111111
$ws.Range("R2").Value = 0.021079197883
$ws.Range("S2").Value = 0.0001687096449493871
$ws.Range("T2").Value = 0.0001687096449493871

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cfh"
$ws.Range("C3").Value = "Itgam"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.4128076666666667
$ws.Range("H3").Value = 1.238423
$ws.Range("I3").Value = 0.001366259689176221
$ws.Range("J3").Value = 0.001366259689176221
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.04027333333333333
$ws.Range("N3").Value = 0.12082
$ws.Range("O3").Value = 0.8765171465674219
$ws.Range("P3").Value = 0.876517146567422
$ws.Range("Q3").Value = 0.01662514076222222
$ws.Range("R3").Value = 0.14962626686
$ws.Range("S3").Value = 0.001197550044226834
$ws.Range("T3").Value = 0.001197550044226834

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cfh"
$ws.Range("C4").Value = "Itgam"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 274.2518313333333
$ws.Range("H4").Value = 822.755494
$ws.Range("I4").Value = 0.90768474543873
$ws.Range("J4").Value = 0.9076847454387301
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.005673666666666667
$ws.Range("N4").Value = 0.017021
$ws.Range("O4").Value = 0.1234828534325781
$ws.Range("P4").Value = 0.1234828534325781
$ws.Range("Q4").Value = 1.556013473708222
$ws.Range("R4").Value = 14.004121263374
$ws.Range("S4").Value = 0.1120835023839977
$ws.Range("T4").Value = 0.1120835023839977

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cfh"
$ws.Range("C5").Value = "Itgam"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 274.2518313333333
$ws.Range("H5").Value = 822.755494
$ws.Range("I5").Value = 0.90768474543873
$ws.Range("J5").Value = 0.9076847454387301
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.04027333333333333
$ws.Range("N5").Value = 0.12082
$ws.Range("O5").Value = 0.8765171465674219
$ws.Range("P5").Value = 0.876517146567422
$ws.Range("Q5").Value = 11.04503542056444
$ws.Range("R5").Value = 99.40531878508
$ws.Range("S5").Value = 0.7956012430547323
$ws.Range("T5").Value = 0.7956012430547325

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Cfh"
$ws.Range("C6").Value = "Itgam"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 27.47972633333333
$ws.Range("H6").Value = 82.439179
$ws.Range("I6").Value = 0.09094899487209368
$ws.Range("J6").Value = 0.09094899487209368
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.005673666666666667
$ws.Range("N6").Value = 0.017021
$ws.Range("O6").Value = 0.1234828534325781
$ws.Range("P6").Value = 0.1234828534325781
$ws.Range("Q6").Value = 0.1559108073065556
$ws.Range("R6").Value = 1.403197265759
$ws.Range("S6").Value = 0.01123064140363104
$ws.Range("T6").Value = 0.01123064140363104

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Cfh"
$ws.Range("C7").Value = "Itgam"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 27.47972633333333
$ws.Range("H7").Value = 82.439179
$ws.Range("I7").Value = 0.09094899487209368
$ws.Range("J7").Value = 0.09094899487209368
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.04027333333333333
$ws.Range("N7").Value = 0.12082
$ws.Range("O7").Value = 0.8765171465674219
$ws.Range("P7").Value = 0.876517146567422
$ws.Range("Q7").Value = 1.106700178531111
$ws.Range("R7").Value = 9.96030160678
$ws.Range("S7").Value = 0.07971835346846264
$ws.Range("T7").Value = 0.07971835346846265

# Remove obsolete rows that held the old "ECs" target-cluster records (old rows 8-10)
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(8).Delete()
